# Insert one new data row into the "Hortaliza, Terminal La Palmera de La
# Serena - Papa" sheet. The new record is placed at row 350 (pushing the
# previous rows 350-409 down to 351-410), growing the used range from
# A1:R409 to A1:R410.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 350..409 down to 351..410 by inserting a blank row.
$ws.Rows.Item(350).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(350, 1).Value  = 8
$ws.Cells.Item(350, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(350, 3).Value  = "Coquimbo"
$ws.Cells.Item(350, 4).Value  = 44694
$ws.Cells.Item(350, 5).Value  = 4
$ws.Cells.Item(350, 6).Value  = 100114001
$ws.Cells.Item(350, 7).Value  = "Papa"
$ws.Cells.Item(350, 8).Value  = "Asterix"
$ws.Cells.Item(350, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(350, 10).Value = 2560
$ws.Cells.Item(350, 11).Value = 8500
$ws.Cells.Item(350, 12).Value = 9000
$ws.Cells.Item(350, 13).Value = 8750
$ws.Cells.Item(350, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(350, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(350, 16).Value = 350
$ws.Cells.Item(350, 17).Value = 25
$ws.Cells.Item(350, 18).Value = "Hortaliza"
